# Weekly update: a new daily price record (row 216) is inserted into the
# "Pepino ensalada" (Vega Modelo de Temuco) dataset, pushing all subsequent
# rows down by one (old row 216 becomes 217, ..., old row 299 becomes 300).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 216; existing rows 216-299 shift down to 217-300.
$ws.Rows("216:216").Insert()

# Populate the newly inserted row 216 with the new record's data.
$ws.Cells.Item(216, 1).Value  = 10
$ws.Cells.Item(216, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(216, 3).Value  = "La Araucanía"
$ws.Cells.Item(216, 4).Value  = 44468
$ws.Cells.Item(216, 5).Value  = 9
$ws.Cells.Item(216, 6).Value  = 100112043
$ws.Cells.Item(216, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(216, 8).Value  = "Sin especificar"
$ws.Cells.Item(216, 9).Value  = "Primera"
$ws.Cells.Item(216, 10).Value = 130
$ws.Cells.Item(216, 11).Value = 18000
$ws.Cells.Item(216, 12).Value = 20000
$ws.Cells.Item(216, 13).Value = 18769
$ws.Cells.Item(216, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(216, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(216, 16).Value = 313
$ws.Cells.Item(216, 17).Value = 60
$ws.Cells.Item(216, 18).Value = "Hortaliza"
